$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Fitness" (column C) values for rows 2-124 (Generation 0-122 of Run 7),
# as captured from the log's best-fitness-so-far re-run.
$updates = @{
    2 = 12316
    3 = 12316
    4 = 12316
    5 = 12224
    6 = 12224
    7 = 11095
    8 = 11095
    9 = 11095
    10 = 10883
    11 = 10883
    12 = 10490
    13 = 10490
    14 = 10445
    15 = 10205
    16 = 9826
    17 = 9826
    18 = 9557
    19 = 9489
    20 = 9489
    21 = 9489
    22 = 9489
    23 = 9489
    24 = 9489
    25 = 9489
    26 = 9489
    27 = 9152
    28 = 9152
    29 = 9152
    30 = 8939
    31 = 8939
    32 = 8939
    33 = 8939
    34 = 8939
    35 = 8939
    36 = 8844
    37 = 8844
    38 = 8517
    39 = 8517
    40 = 8517
    41 = 8517
    42 = 8504
    43 = 8382
    44 = 8012
    45 = 8012
    46 = 8012
    47 = 8012
    48 = 8012
    49 = 8012
    50 = 8012
    51 = 8012
    52 = 8012
    53 = 8012
    54 = 8012
    55 = 8012
    56 = 8012
    57 = 8012
    58 = 8012
    59 = 8012
    60 = 8012
    61 = 8010
    62 = 8010
    63 = 8010
    64 = 8010
    65 = 8010
    66 = 8010
    67 = 8010
    68 = 7970
    69 = 7892
    70 = 7892
    71 = 7892
    72 = 7892
    73 = 7892
    74 = 7892
    75 = 7892
    76 = 7892
    77 = 7892
    78 = 7892
    79 = 7892
    80 = 7892
    81 = 7892
    82 = 7892
    83 = 7892
    84 = 7892
    85 = 7892
    86 = 7892
    87 = 7622
    88 = 7622
    89 = 7622
    90 = 7622
    91 = 7622
    92 = 7622
    93 = 7622
    94 = 7622
    95 = 7622
    96 = 7622
    97 = 7622
    98 = 7622
    99 = 7622
    100 = 7622
    101 = 7622
    102 = 7622
    103 = 7622
    104 = 7622
    105 = 7622
    106 = 7622
    107 = 7622
    108 = 7622
    109 = 7622
    110 = 7622
    111 = 7622
    112 = 7622
    113 = 7622
    114 = 7622
    115 = 7622
    116 = 7622
    117 = 7622
    118 = 7622
    119 = 7622
    120 = 7622
    121 = 7622
    122 = 7573
    123 = 7573
    124 = 7573
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}
